# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (F column) numbers on all four sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 845
$ws.Cells.Item(3, 6).Value = 571
$ws.Cells.Item(6, 6).Value = 1167
$ws.Cells.Item(8, 6).Value = 61
$ws.Cells.Item(11, 6).Value = 1234
$ws.Cells.Item(14, 6).Value = 919
$ws.Cells.Item(15, 6).Value = 898
$ws.Cells.Item(16, 6).Value = 188
$ws.Cells.Item(17, 6).Value = 80
$ws.Cells.Item(18, 6).Value = 84
$ws.Cells.Item(20, 6).Value = 833
$ws.Cells.Item(21, 6).Value = 1760
$ws.Cells.Item(22, 6).Value = 3256
$ws.Cells.Item(23, 6).Value = 963
$ws.Cells.Item(24, 6).Value = 94
$ws.Cells.Item(25, 6).Value = 2349
$ws.Cells.Item(27, 6).Value = 18
$ws.Cells.Item(28, 6).Value = 3222
$ws.Cells.Item(29, 6).Value = 672
$ws.Cells.Item(30, 6).Value = 812
$ws.Cells.Item(32, 6).Value = 526
$ws.Cells.Item(34, 6).Value = 756
$ws.Cells.Item(35, 6).Value = 153
$ws.Cells.Item(37, 6).Value = 89
$ws.Cells.Item(39, 6).Value = 1158
$ws.Cells.Item(40, 6).Value = 1833
$ws.Cells.Item(41, 6).Value = 430
$ws.Cells.Item(43, 6).Value = 565
$ws.Cells.Item(44, 6).Value = 215
$ws.Cells.Item(45, 6).Value = 143
$ws.Cells.Item(47, 6).Value = 60

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(12, 6).Value = 99

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 146

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 845
$ws.Cells.Item(3, 6).Value = 571
$ws.Cells.Item(5, 6).Value = 1167
$ws.Cells.Item(8, 6).Value = 1234
$ws.Cells.Item(10, 6).Value = 919
$ws.Cells.Item(11, 6).Value = 898
$ws.Cells.Item(14, 6).Value = 80
$ws.Cells.Item(16, 6).Value = 84
$ws.Cells.Item(17, 6).Value = 833
$ws.Cells.Item(18, 6).Value = 1760
$ws.Cells.Item(19, 6).Value = 3256
$ws.Cells.Item(20, 6).Value = 963
$ws.Cells.Item(21, 6).Value = 94
$ws.Cells.Item(23, 6).Value = 2349
$ws.Cells.Item(24, 6).Value = 18
$ws.Cells.Item(25, 6).Value = 3222
$ws.Cells.Item(26, 6).Value = 672
$ws.Cells.Item(27, 6).Value = 812
$ws.Cells.Item(34, 6).Value = 99
$ws.Cells.Item(35, 6).Value = 756
$ws.Cells.Item(36, 6).Value = 153
$ws.Cells.Item(38, 6).Value = 89
$ws.Cells.Item(41, 6).Value = 1158
$ws.Cells.Item(42, 6).Value = 1833
$ws.Cells.Item(44, 6).Value = 430
$ws.Cells.Item(45, 6).Value = 565
$ws.Cells.Item(46, 6).Value = 215
$ws.Cells.Item(47, 6).Value = 143
$ws.Cells.Item(49, 6).Value = 60
